$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 (row 1 is the header) through the last used row.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $jCell = $ws.Cells.Item($r, 10)
    $kCell = $ws.Cells.Item($r, 11)

    $jVal = $jCell.Value()
    if ($null -ne $jVal -and $jVal -ne "") {
        $jCell.Value = "'" + ([string]$jVal).TrimEnd(" ")
    }

    $kVal = $kCell.Value()
    if ($null -ne $kVal -and $kVal -ne "") {
        $kCell.Value = " " + ([string]$kVal)
    }
}
